# Update countries & provincias Spain
#
# Daily refresh of the COVID "Pais" dashboard: totals moved for a handful
# of countries, which re-sorts a few adjacent rows (the sheet is kept in
# descending order by total cases, column B) and refreshes the underlying
# case/recovered/critical/death counts for the affected countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4) : updated totals ---
$ws.Cells.Item(4,2).Value = 2089825
$ws.Cells.Item(4,3).Value = 124
$ws.Cells.Item(4,4).Value = 816174
$ws.Cells.Item(4,5).Value = 1157616
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = 116035

# --- Rusia (row 6) : updated totals ---
$ws.Cells.Item(6,2).Value = 511423
$ws.Cells.Item(6,3).Value = 8987
$ws.Cells.Item(6,4).Value = 269370
$ws.Cells.Item(6,5).Value = 235338
$ws.Cells.Item(6,7).Value = 183
$ws.Cells.Item(6,8).Value = 6715

# --- Egipto / Singapur (rows 32-33) swap rank: Singapur overtakes Egipto ---
$ws.Cells.Item(32,1).Value = "Singapur"
$ws.Cells.Item(32,2).Value = 39850
$ws.Cells.Item(32,3).Value = 463
$ws.Cells.Item(32,4).Value = 27286
$ws.Cells.Item(32,5).Value = 12539
$ws.Cells.Item(32,8).Value = 25

$ws.Cells.Item(33,1).Value = "Egipto"
$ws.Cells.Item(33,2).Value = 39726
$ws.Cells.Item(33,4).Value = 10691
$ws.Cells.Item(33,5).Value = 27658
$ws.Cells.Item(33,8).Value = 1377

# --- Ucrania (row 38) : updated totals ---
$ws.Cells.Item(38,2).Value = 29753
$ws.Cells.Item(38,3).Value = 683
$ws.Cells.Item(38,4).Value = 13567
$ws.Cells.Item(38,5).Value = 15316
$ws.Cells.Item(38,7).Value = 16
$ws.Cells.Item(38,8).Value = 870

# --- Afganistan (row 43) : updated totals ---
$ws.Cells.Item(43,2).Value = 23546
$ws.Cells.Item(43,3).Value = 656
$ws.Cells.Item(43,4).Value = 3928
$ws.Cells.Item(43,5).Value = 19172

# --- Austria / Irak / Barein (rows 50-52) : Barein overtakes Austria and Irak ---
$ws.Cells.Item(50,1).Value = "Barein"
$ws.Cells.Item(50,2).Value = 17269
$ws.Cells.Item(50,3).Value = 602
$ws.Cells.Item(50,4).Value = 11903
$ws.Cells.Item(50,5).Value = 5330
$ws.Cells.Item(50,7).Value = 2
$ws.Cells.Item(50,8).Value = 36

$ws.Cells.Item(51,1).Value = "Austria"
$ws.Cells.Item(51,2).Value = 17034
$ws.Cells.Item(51,4).Value = 15949
$ws.Cells.Item(51,5).Value = 411
$ws.Cells.Item(51,8).Value = 674

$ws.Cells.Item(52,1).Value = "Irak"
$ws.Cells.Item(52,2).Value = 16675
$ws.Cells.Item(52,4).Value = 6568
$ws.Cells.Item(52,5).Value = 9650
$ws.Cells.Item(52,8).Value = 457

# --- Armenia (row 54) : updated totals ---
$ws.Cells.Item(54,2).Value = 15281
$ws.Cells.Item(54,3).Value = 612
$ws.Cells.Item(54,4).Value = 5639
$ws.Cells.Item(54,5).Value = 9384
$ws.Cells.Item(54,7).Value = 13
$ws.Cells.Item(54,8).Value = 258

# --- Luxemburgo / Hungria (rows 82-83) swap rank: Hungria overtakes Luxemburgo ---
$ws.Cells.Item(82,1).Value = "Hungria"
$ws.Cells.Item(82,2).Value = 4053
$ws.Cells.Item(82,3).Value = 14
$ws.Cells.Item(82,4).Value = 2447
$ws.Cells.Item(82,5).Value = 1051
$ws.Cells.Item(82,7).Value = 2
$ws.Cells.Item(82,8).Value = 555

$ws.Cells.Item(83,1).Value = "Luxemburgo"
$ws.Cells.Item(83,2).Value = 4052
$ws.Cells.Item(83,4).Value = 3910
$ws.Cells.Item(83,5).Value = 32
$ws.Cells.Item(83,8).Value = 110

# --- Estonia (row 101) : updated totals ---
$ws.Cells.Item(101,2).Value = 1970
$ws.Cells.Item(101,3).Value = 5
$ws.Cells.Item(101,4).Value = 1703
$ws.Cells.Item(101,5).Value = 198

# --- Lituania (row 105) : updated totals ---
$ws.Cells.Item(105,2).Value = 1756
$ws.Cells.Item(105,3).Value = 4
$ws.Cells.Item(105,4).Value = 1400
$ws.Cells.Item(105,5).Value = 282

# --- Eslovaquia (row 108) : updated totals ---
$ws.Cells.Item(108,2).Value = 1542
$ws.Cells.Item(108,3).Value = 1
$ws.Cells.Item(108,5).Value = 105

# --- Letonia (row 122) : updated totals ---
$ws.Cells.Item(122,2).Value = 1096
$ws.Cells.Item(122,3).Value = 2
$ws.Cells.Item(122,5).Value = 251
$ws.Cells.Item(122,7).Value = 1
$ws.Cells.Item(122,8).Value = 27

# --- Georgia (row 133) : updated totals ---
$ws.Cells.Item(133,2).Value = 837
$ws.Cells.Item(133,3).Value = 6
$ws.Cells.Item(133,4).Value = 697
$ws.Cells.Item(133,5).Value = 127

# --- Santo Tome y Principe (row 140) : updated totals ---
$ws.Cells.Item(140,2).Value = 639
$ws.Cells.Item(140,3).Value = 7
$ws.Cells.Item(140,4).Value = 156
$ws.Cells.Item(140,5).Value = 471
